# Update the example covariates data (age, sex, edu) for the connectivity
# example dataset (GroupName2), row by row, leaving the ID column untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 53
$ws.Range("C2").Value = "Female"
$ws.Range("D2").Value = 10
$ws.Range("B3").Value = 60
$ws.Range("C3").Value = "Female"
$ws.Range("D3").Value = 8
$ws.Range("B4").Value = 63
$ws.Range("C4").Value = "Female"
$ws.Range("D4").Value = 13
$ws.Range("B5").Value = 61
$ws.Range("C5").Value = "Male"
$ws.Range("D5").Value = 16
$ws.Range("B6").Value = 70
$ws.Range("C6").Value = "Female"
$ws.Range("D6").Value = 7
$ws.Range("B7").Value = 68
$ws.Range("C7").Value = "Male"
$ws.Range("D7").Value = 14
$ws.Range("B8").Value = 57
$ws.Range("C8").Value = "Male"
$ws.Range("D8").Value = 10
$ws.Range("B9").Value = 57
$ws.Range("C9").Value = "Male"
$ws.Range("D9").Value = 9
$ws.Range("B10").Value = 74
$ws.Range("C10").Value = "Male"
$ws.Range("D10").Value = 11
$ws.Range("B11").Value = 54
$ws.Range("C11").Value = "Female"
$ws.Range("D11").Value = 9
$ws.Range("B12").Value = 63
$ws.Range("C12").Value = "Female"
$ws.Range("D12").Value = 9
$ws.Range("B13").Value = 67
$ws.Range("C13").Value = "Male"
$ws.Range("D13").Value = 6
$ws.Range("B14").Value = 65
$ws.Range("C14").Value = "Female"
$ws.Range("D14").Value = 21
$ws.Range("B15").Value = 79
$ws.Range("C15").Value = "Male"
$ws.Range("D15").Value = 12
$ws.Range("B16").Value = 55
$ws.Range("C16").Value = "Male"
$ws.Range("D16").Value = 12
$ws.Range("B17").Value = 64
$ws.Range("C17").Value = "Male"
$ws.Range("D17").Value = 11
$ws.Range("B18").Value = 76
$ws.Range("C18").Value = "Female"
$ws.Range("D18").Value = 9
$ws.Range("B19").Value = 63
$ws.Range("C19").Value = "Male"
$ws.Range("D19").Value = 16
$ws.Range("B20").Value = 55
$ws.Range("C20").Value = "Female"
$ws.Range("D20").Value = 8
$ws.Range("B21").Value = 61
$ws.Range("C21").Value = "Female"
$ws.Range("D21").Value = 10
$ws.Range("B22").Value = 79
$ws.Range("C22").Value = "Female"
$ws.Range("D22").Value = 15
$ws.Range("B23").Value = 60
$ws.Range("C23").Value = "Female"
$ws.Range("D23").Value = 11
$ws.Range("B24").Value = 76
$ws.Range("C24").Value = "Female"
$ws.Range("D24").Value = 9
$ws.Range("B25").Value = 78
$ws.Range("C25").Value = "Female"
$ws.Range("D25").Value = 7
$ws.Range("B26").Value = 74
$ws.Range("C26").Value = "Female"
$ws.Range("D26").Value = 10
$ws.Range("B27").Value = 58
$ws.Range("C27").Value = "Male"
$ws.Range("D27").Value = 16
$ws.Range("B28").Value = 74
$ws.Range("C28").Value = "Female"
$ws.Range("D28").Value = 14
$ws.Range("B29").Value = 70
$ws.Range("C29").Value = "Female"
$ws.Range("D29").Value = 10
$ws.Range("B30").Value = 55
$ws.Range("C30").Value = "Male"
$ws.Range("D30").Value = 12
$ws.Range("B31").Value = 79
$ws.Range("C31").Value = "Female"
$ws.Range("D31").Value = 10
$ws.Range("B32").Value = 68
$ws.Range("C32").Value = "Female"
$ws.Range("D32").Value = 8
$ws.Range("B33").Value = 74
$ws.Range("C33").Value = "Female"
$ws.Range("D33").Value = 11
$ws.Range("B34").Value = 68
$ws.Range("C34").Value = "Male"
$ws.Range("D34").Value = 11
$ws.Range("B35").Value = 56
$ws.Range("C35").Value = "Male"
$ws.Range("D35").Value = 13
$ws.Range("B36").Value = 78
$ws.Range("C36").Value = "Female"
$ws.Range("D36").Value = 6
$ws.Range("B37").Value = 58
$ws.Range("C37").Value = "Male"
$ws.Range("D37").Value = 6
$ws.Range("B38").Value = 51
$ws.Range("C38").Value = "Male"
$ws.Range("D38").Value = 8
$ws.Range("B39").Value = 59
$ws.Range("C39").Value = "Male"
$ws.Range("D39").Value = 13
$ws.Range("B40").Value = 77
$ws.Range("C40").Value = "Male"
$ws.Range("D40").Value = 7
$ws.Range("B41").Value = 79
$ws.Range("C41").Value = "Female"
$ws.Range("D41").Value = 6
$ws.Range("B42").Value = 77
$ws.Range("C42").Value = "Male"
$ws.Range("D42").Value = 10
$ws.Range("B43").Value = 80
$ws.Range("C43").Value = "Male"
$ws.Range("D43").Value = 8
$ws.Range("B44").Value = 78
$ws.Range("C44").Value = "Male"
$ws.Range("D44").Value = 15
$ws.Range("B45").Value = 65
$ws.Range("C45").Value = "Male"
$ws.Range("D45").Value = 9
$ws.Range("B46").Value = 71
$ws.Range("C46").Value = "Female"
$ws.Range("D46").Value = 10
$ws.Range("B47").Value = 59
$ws.Range("C47").Value = "Female"
$ws.Range("D47").Value = 8
$ws.Range("B48").Value = 73
$ws.Range("C48").Value = "Female"
$ws.Range("D48").Value = 8
$ws.Range("B49").Value = 55
$ws.Range("C49").Value = "Female"
$ws.Range("D49").Value = 6
$ws.Range("B50").Value = 57
$ws.Range("C50").Value = "Female"
$ws.Range("D50").Value = 6
$ws.Range("B51").Value = 63
$ws.Range("C51").Value = "Male"
$ws.Range("D51").Value = 8
